$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: fill in G12/H12 (existing style s=4 stays the same) ---
$ws.Range("G12").Value2 = 9918.8912
$ws.Range("H12").Value2 = "10k"

# --- Row 14: fill in G14/H14. G14 needs a brand new style (new fill, theme
#     "Background1"/white, centered) while H14 keeps the existing s=4 style ---
$ws.Range("H14").Value2 = "10k"

$ws.Range("G14").Value2 = 9918
$ws.Range("G14").Interior.ThemeColor = 2
$ws.Range("G14").Interior.TintAndShade = 0

# --- Row 15: update existing value ---
$ws.Range("G15").Value2 = 6408

# --- New rows 20 and 21 with 9x9 "10k" measurement averages ---

# Row 20
$ws.Range("A11").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value2 = "single"
$ws.Range("B20").Value2 = "no"
$ws.Range("C20").Value2 = "make path free"
$ws.Range("G20").Value2 = 70278
$ws.Range("H20").Value2 = "10k"

# Row 21
$ws.Range("A12").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value2 = "max"
$ws.Range("B21").Value2 = "no"
$ws.Range("C21").Value2 = "make path free"
$ws.Range("G21").Value2 = 8849
$ws.Range("H21").Value2 = "10k"

$excel.CutCopyMode = 0

# --- Update selection to match the final state ---
$null = $ws.Range("F21").Select()
